# Update "想去人数" (want-to-go count) figures across the four sheets to
# reflect the newly scraped data (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# ---- Sheet 1: 展览 (Exhibitions) ----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(4, 6).Value  = 390
$ws1.Cells.Item(6, 6).Value  = 6
$ws1.Cells.Item(7, 6).Value  = 557
$ws1.Cells.Item(8, 6).Value  = 54
$ws1.Cells.Item(9, 6).Value  = 9801
$ws1.Cells.Item(11, 6).Value = 2653
$ws1.Cells.Item(13, 6).Value = 2391
$ws1.Cells.Item(14, 6).Value = 2660
$ws1.Cells.Item(16, 6).Value = 276
$ws1.Cells.Item(17, 6).Value = 2085
$ws1.Cells.Item(22, 6).Value = 72
$ws1.Cells.Item(23, 6).Value = 300
$ws1.Cells.Item(25, 6).Value = 148
$ws1.Cells.Item(27, 6).Value = 1290
$ws1.Cells.Item(28, 6).Value = 1241
$ws1.Cells.Item(32, 6).Value = 1683
$ws1.Cells.Item(33, 6).Value = 2817
$ws1.Cells.Item(35, 6).Value = 993
$ws1.Cells.Item(36, 6).Value = 359
$ws1.Cells.Item(38, 6).Value = 1277
$ws1.Cells.Item(39, 6).Value = 53
$ws1.Cells.Item(40, 6).Value = 55
$ws1.Cells.Item(41, 6).Value = 51
$ws1.Cells.Item(42, 6).Value = 28

# ---- Sheet 2: 演出 (Performances) ----
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(2, 6).Value = 8

# ---- Sheet 3: 本地生活 (Local life) ----
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Cells.Item(2, 6).Value = 723
$ws3.Cells.Item(3, 6).Value = 953
$ws3.Cells.Item(5, 6).Value = 1728
$ws3.Cells.Item(5, 7).Value = "暂时售罄"

# ---- Sheet 4: 全部类型 (All types) ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(3, 6).Value  = 724
$ws4.Cells.Item(4, 6).Value  = 953
$ws4.Cells.Item(6, 6).Value  = 390
$ws4.Cells.Item(7, 6).Value  = 8
$ws4.Cells.Item(10, 6).Value = 557
$ws4.Cells.Item(11, 6).Value = 54
$ws4.Cells.Item(12, 6).Value = 9801
$ws4.Cells.Item(16, 6).Value = 2653
$ws4.Cells.Item(18, 6).Value = 2391
$ws4.Cells.Item(19, 6).Value = 2660
$ws4.Cells.Item(20, 6).Value = 276
$ws4.Cells.Item(21, 6).Value = 2085
$ws4.Cells.Item(26, 6).Value = 300
$ws4.Cells.Item(28, 6).Value = 148
$ws4.Cells.Item(30, 6).Value = 1290
$ws4.Cells.Item(31, 6).Value = 1241
$ws4.Cells.Item(35, 6).Value = 1683
$ws4.Cells.Item(37, 6).Value = 2817
$ws4.Cells.Item(38, 6).Value = 993
$ws4.Cells.Item(41, 6).Value = 359
$ws4.Cells.Item(45, 6).Value = 1277
$ws4.Cells.Item(46, 6).Value = 51
$ws4.Cells.Item(47, 6).Value = 28
